# Automatische test-sync: 2025-08-26 20:29:50
#
# Adds the new "Retour / Terugbetaling" mail-log entry (row 6) to the Logs
# sheet, the matching dashboard tally row (row 3) on the Dashboard sheet,
# extends the conditional formatting ranges to cover the new row, and
# extends the bar chart's category/value series to include the new
# Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new log row (row 6)
# ---------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A6").Value = "Retour status"
$wsLogs.Range("B6").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("D6").Value = "Retour / Terugbetaling"
$wsLogs.Range("F6").Value = "2025-08-26 20:29:43"
$wsLogs.Range("G6").Value = "Nee"
$wsLogs.Range("H6").Value = "Ja"
$wsLogs.Range("I6").Value = "Nee"
$wsLogs.Range("J6").Value = "Nee"

# Extend the conditional formatting ranges on the Logs sheet so they keep
# covering column D/G/H/I/J through the newly added row 6.
$cfD = $wsLogs.Range("D2:D5").FormatConditions
for ($i = 1; $i -le $cfD.Count; $i++) {
    $cfD.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D6"))
}

$cfG = $wsLogs.Range("G2:G5").FormatConditions
for ($i = 1; $i -le $cfG.Count; $i++) {
    $cfG.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G6"))
}

$cfH = $wsLogs.Range("H2:H5").FormatConditions
for ($i = 1; $i -le $cfH.Count; $i++) {
    $cfH.Item($i).ModifyAppliesToRange($wsLogs.Range("H2:H6"))
}

$cfI = $wsLogs.Range("I2:I5").FormatConditions
for ($i = 1; $i -le $cfI.Count; $i++) {
    $cfI.Item($i).ModifyAppliesToRange($wsLogs.Range("I2:I6"))
}

$cfJ = $wsLogs.Range("J2:J5").FormatConditions
for ($i = 1; $i -le $cfJ.Count; $i++) {
    $cfJ.Item($i).ModifyAppliesToRange($wsLogs.Range("J2:J6"))
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: append the new tally row (row 3)
# ---------------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A3").Value = "Retour / Terugbetaling"
$wsDash.Range("B3").Value = 1

# ---------------------------------------------------------------------
# 3) Chart: extend the category/value series references to include the
#    new Dashboard row (A2:A3 / B2:B3 instead of single cells A2 / B2).
# ---------------------------------------------------------------------
$chartObjs = $wsDash.ChartObjects()
$chartObj = $chartObjs.Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection()
$ser = $series.Item(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
